$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment each date serial value in F2:F7 by 1 day,
# matching the diff (final edits and rerender).
for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value2 = $cell.Value2 + 1
}
